$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting the existing rows 23:48 down to 24:49
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new record's data
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 44638
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 100112043
$ws.Range("G23").Value = "Pepino dulce"
$ws.Range("H23").Value = "Cultivar IV Región"
$ws.Range("I23").Value = "Especial"
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 21000
$ws.Range("M23").Value = 20500
$ws.Range("N23").Value = "$/bandeja 18 kilos"
$ws.Range("O23").Value = "Provincia de Limarí"
$ws.Range("P23").Value = 1139
$ws.Range("Q23").Value = 18
$ws.Range("R23").Value = "Hortaliza"
